$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.578.50"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "3.406.65"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.402.02"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.72%  "
$ws.Range("E10").Value = "  -9.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.03"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.77%  "
$ws.Range("E12").Value = "  -9.98%  "
$ws.Range("D13").Value = "3.974.71"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("E14").Value = "  -10.64%  "
$ws.Range("D15").Value = "3.425.66"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "64.489.27"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.81"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -10.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -16.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.32"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -11.11%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.535"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -10.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.42"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.88%  "
$ws.Range("D26").Value = "3.536.25"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("E27").Value = "  -12.15%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  -11.19%  "
$ws.Range("E30").Value = "  -12.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.41%  "
$ws.Range("D32").Value = "3.418.51"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.81%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.140"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -10.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.22"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -13.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.59"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -14.45%  "
$ws.Range("E39").Value = "  -13.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.51"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -15.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0753"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -8.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.798"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.93"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.84%  "
$ws.Range("E45").Value = "  -16.12%  "
$ws.Range("E46").Value = "  -11.10%  "
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.28%  "
$ws.Range("E49").Value = "  -10.04%  "
$ws.Range("D50").Value = "2.187.62"
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -19.25%  "
